# Rename/position handled via adding a new "metadata" worksheet after "data",
# updating the time_taken (F column) timestamps on "data", and populating
# "metadata" with the panel query info (mirrors the commit's refined export).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1) Refresh the "time_taken" timestamps in column F (rows 2-120) of the
#    "data" sheet.
# ---------------------------------------------------------------------
$newTimestamps = @(
    "2021-10-05 14:21:06.412627",
    "2021-10-05 14:21:06.412637",
    "2021-10-05 14:21:06.412640",
    "2021-10-05 14:21:06.412643",
    "2021-10-05 14:21:06.412646",
    "2021-10-05 14:21:06.412650",
    "2021-10-05 14:21:06.412653",
    "2021-10-05 14:21:06.412656",
    "2021-10-05 14:21:06.412659",
    "2021-10-05 14:21:06.412662",
    "2021-10-05 14:21:06.412665",
    "2021-10-05 14:21:06.412667",
    "2021-10-05 14:21:06.412670",
    "2021-10-05 14:21:06.412673",
    "2021-10-05 14:21:06.412675",
    "2021-10-05 14:21:06.412678",
    "2021-10-05 14:21:06.412681",
    "2021-10-05 14:21:06.412684",
    "2021-10-05 14:21:06.412687",
    "2021-10-05 14:21:06.412690",
    "2021-10-05 14:21:06.412692",
    "2021-10-05 14:21:06.412695",
    "2021-10-05 14:21:06.412698",
    "2021-10-05 14:21:06.412701",
    "2021-10-05 14:21:06.412704",
    "2021-10-05 14:21:06.412707",
    "2021-10-05 14:21:06.412710",
    "2021-10-05 14:21:06.412713",
    "2021-10-05 14:21:06.412716",
    "2021-10-05 14:21:06.412718",
    "2021-10-05 14:21:06.412721",
    "2021-10-05 14:21:06.412724",
    "2021-10-05 14:21:06.412727",
    "2021-10-05 14:21:06.412730",
    "2021-10-05 14:21:06.412733",
    "2021-10-05 14:21:06.412735",
    "2021-10-05 14:21:06.412738",
    "2021-10-05 14:21:06.412741",
    "2021-10-05 14:21:06.412744",
    "2021-10-05 14:21:06.412747",
    "2021-10-05 14:21:06.412750",
    "2021-10-05 14:21:06.412753",
    "2021-10-05 14:21:06.412756",
    "2021-10-05 14:21:06.412758",
    "2021-10-05 14:21:06.412761",
    "2021-10-05 14:21:06.412764",
    "2021-10-05 14:21:06.412767",
    "2021-10-05 14:21:06.412770",
    "2021-10-05 14:21:06.412773",
    "2021-10-05 14:21:06.412775",
    "2021-10-05 14:21:06.412778",
    "2021-10-05 14:21:06.412781",
    "2021-10-05 14:21:06.412784",
    "2021-10-05 14:21:06.412787",
    "2021-10-05 14:21:06.412790",
    "2021-10-05 14:21:06.412793",
    "2021-10-05 14:21:06.412795",
    "2021-10-05 14:21:06.412798",
    "2021-10-05 14:21:06.412801",
    "2021-10-05 14:21:06.412803",
    "2021-10-05 14:21:06.412806",
    "2021-10-05 14:21:06.412809",
    "2021-10-05 14:21:06.412811",
    "2021-10-05 14:21:06.412814",
    "2021-10-05 14:21:06.412818",
    "2021-10-05 14:21:06.412821",
    "2021-10-05 14:21:06.412824",
    "2021-10-05 14:21:06.412826",
    "2021-10-05 14:21:06.412829",
    "2021-10-05 14:21:06.412832",
    "2021-10-05 14:21:06.412835",
    "2021-10-05 14:21:06.412837",
    "2021-10-05 14:21:06.412840",
    "2021-10-05 14:21:06.412843",
    "2021-10-05 14:21:06.412846",
    "2021-10-05 14:21:06.412848",
    "2021-10-05 14:21:06.412853",
    "2021-10-05 14:21:06.412856",
    "2021-10-05 14:21:06.412859",
    "2021-10-05 14:21:06.412862",
    "2021-10-05 14:21:06.412865",
    "2021-10-05 14:21:06.412868",
    "2021-10-05 14:21:06.412870",
    "2021-10-05 14:21:06.412873",
    "2021-10-05 14:21:06.412876",
    "2021-10-05 14:21:06.412879",
    "2021-10-05 14:21:06.412882",
    "2021-10-05 14:21:06.412884",
    "2021-10-05 14:21:06.412887",
    "2021-10-05 14:21:06.412890",
    "2021-10-05 14:21:06.412893",
    "2021-10-05 14:21:06.412896",
    "2021-10-05 14:21:06.412899",
    "2021-10-05 14:21:06.412903",
    "2021-10-05 14:21:06.412906",
    "2021-10-05 14:21:06.412908",
    "2021-10-05 14:21:06.412911",
    "2021-10-05 14:21:06.412914",
    "2021-10-05 14:21:06.412917",
    "2021-10-05 14:21:06.412919",
    "2021-10-05 14:21:06.412922",
    "2021-10-05 14:21:06.412925",
    "2021-10-05 14:21:06.412927",
    "2021-10-05 14:21:06.412930",
    "2021-10-05 14:21:06.412933",
    "2021-10-05 14:21:06.412936",
    "2021-10-05 14:21:06.412938",
    "2021-10-05 14:21:06.412941",
    "2021-10-05 14:21:06.412946",
    "2021-10-05 14:21:06.412949",
    "2021-10-05 14:21:06.412952",
    "2021-10-05 14:21:06.412954",
    "2021-10-05 14:21:06.412957",
    "2021-10-05 14:21:06.412960",
    "2021-10-05 14:21:06.412963",
    "2021-10-05 14:21:06.412966",
    "2021-10-05 14:21:06.412969",
    "2021-10-05 14:21:06.412972",
    "2021-10-05 14:21:06.412974"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------------
# 2) Add the new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Inherited bleeding disorders"
$meta.Range("C2").Value = 175

# "1.164" must stay textual (not be coerced to the number 1.164) -- write
# it as a formula that evaluates to a text value, then collapse it down to
# a literal value via copy / paste-values so no formula or stray number
# format survives into the saved file.
$meta.Range("D2").Formula = "=""1.164"""
$meta.Range("D2").Copy()
$meta.Range("D2").PasteSpecial(-4163)

$meta.Range("E2").Value = "2021-08-02T07:43:09.605718Z"
$meta.Range("F2").Value = "2021-10-05 14:21:06.408758"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/175/?format=json"

# ---------------------------------------------------------------------
# 3) Mirror the bold / centred / bordered header style used by "data" onto
#    the new sheet's header row (B1:G1) and its styled A2 index cell.
# ---------------------------------------------------------------------
$data.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("A1").Select()
